# Apply updates to both the "展览" and "全部类型" sheets, which contain
# identical data in this workbook. Update the Cover image URL on row 2,
# and bump the "想去人数" counts on rows 4 and 9.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202406/hZd8mGjR1718691261895.jpeg"
    $ws.Range("F4").Value = 1442
    $ws.Range("F9").Value = 230
}
